$d = $word.ActiveDocument

# Locate the final top-level (i.e. not inside a table) paragraph whose text
# is just "}" -- this is the paragraph right before the two blank paragraphs
# that need to be removed.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "}") {
        $withinTable = $p.Range.Information(12)
        if ($withinTable -eq $false) {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -gt 0) {
    # The two empty paragraphs immediately following the target "}" paragraph
    # are deleted. Deleting the same (now-next) index twice removes both,
    # since the collection re-indexes after each deletion.
    $d.Paragraphs($targetIndex + 1).Range.Delete()
    $d.Paragraphs($targetIndex + 1).Range.Delete()
}
